$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.331414580345154
$ws.Range("B1").Value = 1.532049417495728
$ws.Range("C1").Value = 1.940873265266418
$ws.Range("D1").Value = 2.714714050292969
$ws.Range("E1").Value = 15
